# #3456 updated PM Property ID
# Updates the "Portfolio Manager Building ID" values (column B) on the
# "BPS Data" sheet for rows 2-10, and restores the active selection to
# B2:B10 (active cell B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

$updates = @{
    2  = 22178843
    3  = 22178844
    4  = 22178845
    5  = 22178846
    6  = 22178847
    7  = 22178848
    8  = 22178849
    9  = 22178850
    10 = 22178851
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

$ws.Activate()
$ws.Range("B2:B10").Select()
